$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.298.09'
$ws.Range("E2").Value = '  +4.34%  '
$ws.Range("D3").Value = '2.347.72'
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''546.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.40%  '
$ws.Range("D6").Value = '''132.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.40%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +1.57%  '
$ws.Range("D9").Value = '2.346.44'
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("D11").Value = '''5.53'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("D13").Value = '''0.335'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").Value = '''23.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").Value = '2.763.62'
$ws.Range("E15").Value = '  +3.09%  '
$ws.Range("D16").Value = '60.216.87'
$ws.Range("E16").Value = '  +4.24%  '
$ws.Range("E17").Value = '  +2.25%  '
$ws.Range("D18").Value = '2.354.48'
$ws.Range("E18").Value = '  +4.19%  '
$ws.Range("D19").Value = '''10.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.33%  '
$ws.Range("D20").Value = '''4.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").Value = '''6.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.65%  '
$ws.Range("D22").Value = '''314.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '''63.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").Value = '''0.173'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '''7.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").Value = '''1.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.42%  '
$ws.Range("E29").Value = '  +2.95%  '
$ws.Range("D30").Value = '''171.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").Value = '''1.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.19%  '
$ws.Range("D32").Value = '0.0₃0729'
$ws.Range("E32").Value = '  +2.36%  '
$ws.Range("D33").Value = '''5.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.50%  '
$ws.Range("E34").Value = '  +15.79%  '
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("D36").Value = '''18.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '''4.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.54%  '
$ws.Range("D40").Value = '''325.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.65%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''1.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.21%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''38.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").Value = '''141.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").Value = '''3.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.04%  '
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = '''19.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.14%  '
$ws.Range("D47").Value = '''0.0499'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("D48").Value = '''0.564'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("D50").Value = '0.0₆0213'
$ws.Range("E50").Value = '  +18.60%  '
$ws.Range("E51").Value = '  +0.86%  '
